$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.257.33"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.863.53"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'236.64"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.4705"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("D8").Value = "'0.2913"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").Value = "'0.06554"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "'21.88"
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("D11").Value = "'0.07929"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "'98.03"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.869.10"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'5.164"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "'0.6815"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "'266.84"
$ws.Range("E16").Value = "  -5.00%  "
$ws.Range("D17").Value = "30.256.39"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "'13.78"
$ws.Range("E18").Value = "  +8.59%  "
$ws.Range("D19").Value = "'0.9996"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'0.000007427"
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("D21").Value = "2.114.34"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'5.316"
$ws.Range("E22").Value = "  -3.31%  "
$ws.Range("D23").Value = "'0.9999"
$ws.Range("D24").Value = "'6.180"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "'167.46"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").Value = "'9.236"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").Value = "'18.94"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").Value = "'1.956"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").Value = "'0.09856"
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("D31").Value = "'4.376"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "'1.470"
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").Value = "'4.055"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").Value = "'0.04713"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").Value = "'1.132"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").Value = "'0.7038"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").Value = "'2.703"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'0.01880"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "'2.612"
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("D40").Value = "'6.290"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").Value = "'74.23"
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("D42").Value = "'1.952"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").Value = "'0.8463"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "'0.4164"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").Value = "'0.9986"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "'103.39"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").Value = "'7.176"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").Value = "'952.43"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("D49").Value = "'9.261"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "'34.15"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'0.05657"
$ws.Range("E51").Value = "  +0.42%  "
